# Commit message: "include others for sos"
# The "Include Others" column (E) for the "Facings SOS" row (row 2) should
# be changed from "Exclude" to "Include".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional KPIs")
$ws.Range("E2").Value = "Include"
